$d = $word.ActiveDocument

# Locate the "Frameworks to consider" heading paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r", "`a") -eq "Frameworks to consider") {
        $target = $p
    }
}

# Delete everything from the end of that paragraph through the end of the
# document body (i.e. all the old "frameworks to consider" content: the
# hyperlinks, Flask/Django/FastAPI/Bottle lines and the two pictures).
$deleteStart = $target.Range.End
$deleteEnd = $d.Content.End
if ($deleteEnd -gt $deleteStart) {
    $d.Range($deleteStart, $deleteEnd).Delete()
}

# The heading paragraph is now the last paragraph in the document. Replace
# its text with the new FastAPI documentation link.
[void]$target.Range.Find.Execute("Frameworks to consider", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "https://fastapi.tiangolo.com/", 2)

# Turn it from a Heading 1 into a plain paragraph (matching the removal of
# <w:pStyle w:val="Heading1"/> in the target XML).
$target.Style = $d.Styles.Item("Normal")
